$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-28 23:03:27"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-28 23:03:22"
$wsZhCn.Range("K2").Value = "2016-08-28 23:03:52"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-28 23:03:27"
$wsDeDe.Range("K2").Value = "2016-08-28 23:03:59"
